$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the time-range values in column B
$ws.Range("B9").Value = "20:00 - 20:04"
$ws.Range("B10").Value = "20:05 - 20:09"

# Update the active cell selection shown in the sheet view
$ws.Range("B12").Select()
